$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new line "Django" to the Test Data cell for the TC-006 row (F14)
$ws.Range("F14").Value = "Valid Url" + [char]10 + "Valid API Response" + [char]10 + "Test Data" + [char]10 + "Django"

# Reflect the end-user navigation/selection state after the edit: the view
# scrolled so row 9 is the top visible row, and the active cell moved to G14
# (the cell right after the one that was edited).
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("G14").Select()
